$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '74.844.34'
$ws.Range("E2").Value = '  +1.44%  '
$ws.Range("D3").Value = '2.811.68'
$ws.Range("E3").Value = '  +7.17%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '186.98'
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '593.48'
$ws.Range("E6").Value = '  +1.98%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +2.79%  '
$ws.Range("E9").Value = '  -4.79%  '
$ws.Range("D10").Value = '2.810.33'
$ws.Range("E10").Value = '  +7.25%  '
$ws.Range("E11").Value = '  -1.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.370'
$ws.Range("E12").Value = '  +3.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.85'
$ws.Range("E13").Value = '  +1.74%  '
$ws.Range("D14").Value = '3.326.81'
$ws.Range("E14").Value = '  +7.24%  '
$ws.Range("D15").Value = '74.848.41'
$ws.Range("E15").Value = '  +1.63%  '
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.72'
$ws.Range("E17").Value = '  +2.11%  '
$ws.Range("D18").Value = '2.813.71'
$ws.Range("E18").Value = '  +7.30%  '
$ws.Range("E19").Value = '  -1.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.25'
$ws.Range("E20").Value = '  +3.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.78'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("E22").Value = '  -2.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.05'
$ws.Range("E23").Value = '  -0.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.21'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.58'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").Value = '2.957.77'
$ws.Range("E27").Value = '  +7.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.14'
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.67'
$ws.Range("E29").Value = '  +2.76%  '
$ws.Range("E30").Value = '  +10.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.39'
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '511.14'
$ws.Range("E33").Value = '  -2.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.68'
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("E35").Value = '  +2.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.12'
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.86'
$ws.Range("E38").Value = '  +3.65%  '
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.38'
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '185.17'
$ws.Range("E41").Value = '  +14.79%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.337'
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.98'
$ws.Range("E44").Value = '  +1.40%  '
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.20'
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.73'
$ws.Range("E47").Value = '  +2.11%  '
$ws.Range("E48").Value = '  -1.95%  '
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.568'
$ws.Range("E50").Value = '  +7.52%  '
$ws.Range("E51").Value = '  +1.77%  '
